$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all target cells keep their original text formatting
# (columns D and E store numeric-looking / percentage values as text)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.696.92"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +1.07%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.873.43"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.34%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.005"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "332.08"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +3.42%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4718"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +5.97%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3953"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +3.26%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "47.54"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -2.94%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08031"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +2.70%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.022"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "21.81"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +1.71%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.883.40"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +2.18%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.955"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +2.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.144"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.90%  "
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.51%  "
$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "87.19"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001044"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +2.01%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06630"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.17"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.15%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.004"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.719.49"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.495"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.65%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.01"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.42%  "
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.105.30"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +2.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.34"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +3.16%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.18"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.75%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.098"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.563"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +2.02%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "122.43"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +2.67%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.9673"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +4.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.09567"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +3.04%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.452"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.52%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.628"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.94%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.299"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +1.86%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06115"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +3.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02260"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +2.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.224"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +1.32%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.168"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -1.36%  "
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +0.12%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5982"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +1.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1904"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +3.35%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "10.26"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +0.66%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.250"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.30%  "
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5685"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.98%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "12.23"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.406"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.61%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.929"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.84%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06816"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00000000313"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +8.20%  "
